# Hortaliza / Betarraga - Terminal La Palmera de La Serena
# Weekly update: insert a new week's data (rows 116-117) at the top of the
# historical block, pushing every existing week down by one row-pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 116:273 down to 118:275, inserting two blank rows at 116:117.
$ws.Rows("116:117").Insert()

# Populate the new week's "Primera" row (116).
$ws.Range("A116").Value = 8
$ws.Range("B116").Value = "Terminal La Palmera de La Serena"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44665
$ws.Range("E116").Value = 4
$ws.Range("F116").Value = 100114014
$ws.Range("G116").Value = "Betarraga"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 450
$ws.Range("L116").Value = 500
$ws.Range("M116").Value = 475
$ws.Range("N116").Value = "$/paquete 3 unidades"
$ws.Range("O116").Value = "Provincia del Elquí"
$ws.Range("P116").Value = 158
$ws.Range("Q116").Value = 3
$ws.Range("R116").Value = "Hortaliza"

# Populate the new week's "Segunda" row (117).
$ws.Range("A117").Value = 8
$ws.Range("B117").Value = "Terminal La Palmera de La Serena"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44665
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100114014
$ws.Range("G117").Value = "Betarraga"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Segunda"
$ws.Range("J117").Value = 1500
$ws.Range("K117").Value = 350
$ws.Range("L117").Value = 400
$ws.Range("M117").Value = 375
$ws.Range("N117").Value = "$/paquete 3 unidades"
$ws.Range("O117").Value = "Provincia del Elquí"
$ws.Range("P117").Value = 125
$ws.Range("Q117").Value = 3
$ws.Range("R117").Value = "Hortaliza"
